$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: all financial data cells removed (kept only A2:C2) ---
$ws.Range("D2:AJ2").ClearContents()

# --- Row 3: selected cells removed, rest updated to new values ---
$ws.Range("D3").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("T3").ClearContents()
$ws.Range("W3:Z3").ClearContents()
$ws.Range("AD3").ClearContents()
$ws.Range("AH3").ClearContents()

# Row 3
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 2605
$ws.Range("L3").Value = 1975
$ws.Range("M3").Value = 630
$ws.Range("N3").Value = 628
$ws.Range("O3").Value = 2
$ws.Range("P3").Value = 104
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 47
$ws.Range("S3").Value = 0
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 1264
$ws.Range("AA3").Value = 313.22
$ws.Range("AB3").Value = 505.95
$ws.Range("AC3").Value = -1
$ws.Range("AE3").Value = 1515
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 41458958

# Row 4
$ws.Range("D4").Value = 6402
$ws.Range("E4").Value = 438
$ws.Range("F4").Value = 438
$ws.Range("G4").Value = 430
$ws.Range("H4").Value = 349
$ws.Range("I4").Value = 347
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 4367
$ws.Range("L4").Value = 2439
$ws.Range("M4").Value = 1928
$ws.Range("N4").Value = 1923
$ws.Range("O4").Value = 5
$ws.Range("P4").Value = 135
$ws.Range("Q4").Value = 327
$ws.Range("R4").Value = -552
$ws.Range("S4").Value = 943
$ws.Range("T4").Value = 561
$ws.Range("U4").Value = -234
$ws.Range("V4").Value = 1346
$ws.Range("W4").Value = 6.85
$ws.Range("X4").Value = 5.46
$ws.Range("Y4").Value = 27.19
$ws.Range("Z4").Value = 10.02
$ws.Range("AA4").Value = 126.51
$ws.Range("AB4").Value = 1294.29
$ws.Range("AC4").Value = 775
$ws.Range("AD4").Value = 9.869999999999999
$ws.Range("AE4").Value = 3570
$ws.Range("AF4").Value = 2.14
$ws.Range("AG4").Value = 18
$ws.Range("AH4").Value = 0.23
$ws.Range("AI4").Value = 2.72
$ws.Range("AJ4").Value = 53858958

# Row 5
$ws.Range("D5").Value = 7682
$ws.Range("E5").Value = 585
$ws.Range("F5").Value = 585
$ws.Range("G5").Value = 528
$ws.Range("H5").Value = 416
$ws.Range("I5").Value = 415
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 6268
$ws.Range("L5").Value = 3653
$ws.Range("M5").Value = 2615
$ws.Range("N5").Value = 2580
$ws.Range("O5").Value = 35
$ws.Range("P5").Value = 151
$ws.Range("Q5").Value = 322
$ws.Range("R5").Value = -968
$ws.Range("S5").Value = 341
$ws.Range("T5").Value = 1113
$ws.Range("U5").Value = -791
$ws.Range("V5").Value = 2131
$ws.Range("W5").Value = 7.62
$ws.Range("X5").Value = 5.42
$ws.Range("Y5").Value = 18.41
$ws.Range("Z5").Value = 7.83
$ws.Range("AA5").Value = 139.71
$ws.Range("AB5").Value = 1699.13
$ws.Range("AC5").Value = 750
$ws.Range("AD5").Value = 15.99
$ws.Range("AE5").Value = 4262
$ws.Range("AF5").Value = 2.82
$ws.Range("AG5").Value = 18
$ws.Range("AH5").Value = 0.15
$ws.Range("AI5").Value = 2.56
$ws.Range("AJ5").Value = 60536670

# Row 6
$ws.Range("D6").Value = 8792
$ws.Range("E6").Value = 403
$ws.Range("F6").Value = 403
$ws.Range("G6").Value = 231
$ws.Range("H6").Value = 171
$ws.Range("I6").Value = 172
$ws.Range("K6").Value = 7605
$ws.Range("L6").Value = 4773
$ws.Range("M6").Value = 2831
$ws.Range("N6").Value = 2798
$ws.Range("P6").Value = 303
$ws.Range("Q6").Value = 429
$ws.Range("R6").Value = -1045
$ws.Range("S6").Value = 741
$ws.Range("T6").Value = 1098
$ws.Range("U6").Value = -668
$ws.Range("V6").Value = 3039
$ws.Range("W6").Value = 4.58
$ws.Range("X6").Value = 1.94
$ws.Range("Y6").Value = 6.39
$ws.Range("Z6").Value = 2.46
$ws.Range("AA6").Value = 168.59
$ws.Range("AB6").Value = 854.29
$ws.Range("AC6").Value = 284
$ws.Range("AD6").Value = 37.53
$ws.Range("AE6").Value = 4621
$ws.Range("AF6").Value = 2.3
$ws.Range("AG6").Value = 35
$ws.Range("AH6").Value = 0.33
$ws.Range("AI6").Value = 12.33
$ws.Range("AJ6").Value = 60536670

# Row 7
$ws.Range("D7").Value = 11408
$ws.Range("E7").Value = 805
$ws.Range("G7").Value = 729
$ws.Range("H7").Value = 581
$ws.Range("I7").Value = 576
$ws.Range("K7").Value = 8719
$ws.Range("L7").Value = 5222
$ws.Range("M7").Value = 3499
$ws.Range("N7").Value = 3426
$ws.Range("P7").Value = 301
$ws.Range("Q7").Value = 926
$ws.Range("R7").Value = -632
$ws.Range("S7").Value = -105
$ws.Range("T7").Value = 566
$ws.Range("U7").Value = 32
$ws.Range("W7").Value = 7.05
$ws.Range("X7").Value = 5.09
$ws.Range("Y7").Value = 18.51
$ws.Range("Z7").Value = 7.12
$ws.Range("AA7").Value = 149.25
$ws.Range("AC7").Value = 951
$ws.Range("AD7").Value = 16.61
$ws.Range("AE7").Value = 5660
$ws.Range("AF7").Value = 2.79
$ws.Range("AG7").Value = 41
$ws.Range("AH7").Value = 0.26
$ws.Range("AI7").Value = 4.32

# Row 8
$ws.Range("D8").Value = 13819
$ws.Range("E8").Value = 1121
$ws.Range("G8").Value = 1039
$ws.Range("H8").Value = 796
$ws.Range("I8").Value = 802
$ws.Range("K8").Value = 9595
$ws.Range("L8").Value = 5324
$ws.Range("M8").Value = 4271
$ws.Range("N8").Value = 4218
$ws.Range("P8").Value = 301
$ws.Range("Q8").Value = 1043
$ws.Range("R8").Value = -703
$ws.Range("S8").Value = -154
$ws.Range("T8").Value = 617
$ws.Range("U8").Value = 365
$ws.Range("W8").Value = 8.109999999999999
$ws.Range("X8").Value = 5.76
$ws.Range("Y8").Value = 20.97
$ws.Range("Z8").Value = 8.69
$ws.Range("AA8").Value = 124.66
$ws.Range("AC8").Value = 1324
$ws.Range("AD8").Value = 11.93
$ws.Range("AE8").Value = 6967
$ws.Range("AF8").Value = 2.27
$ws.Range("AG8").Value = 54
$ws.Range("AH8").Value = 0.34
$ws.Range("AI8").Value = 4.11

# Row 9
$ws.Range("D9").Value = 16398
$ws.Range("E9").Value = 1409
$ws.Range("G9").Value = 1342
$ws.Range("H9").Value = 1043
$ws.Range("I9").Value = 1054
$ws.Range("K9").Value = 10820
$ws.Range("L9").Value = 5535
$ws.Range("M9").Value = 5287
$ws.Range("N9").Value = 5252
$ws.Range("P9").Value = 301
$ws.Range("Q9").Value = 1311
$ws.Range("R9").Value = -596
$ws.Range("S9").Value = -257
$ws.Range("T9").Value = 489
$ws.Range("U9").Value = 759
$ws.Range("W9").Value = 8.59
$ws.Range("X9").Value = 6.36
$ws.Range("Y9").Value = 22.26
$ws.Range("Z9").Value = 10.22
$ws.Range("AA9").Value = 104.68
$ws.Range("AC9").Value = 1741
$ws.Range("AD9").Value = 9.07
$ws.Range("AE9").Value = 8676
$ws.Range("AF9").Value = 1.82
$ws.Range("AG9").Value = 67
$ws.Range("AH9").Value = 0.42
$ws.Range("AI9").Value = 3.86
